$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete tail of the log (rows 90-120), since the
# refreshed export only covers up to row 89 (dimension A1:E89).
$ws.Range("A90:E120").EntireRow.Delete() | Out-Null

# Every remaining detail row (2-89) now reflects the "invoice" page view
# that the new Firebase real-time sync produced, only timestamps differ.
$ws.Range("C2:C89").Value = "invoice"
$ws.Range("D2:D89").Value = "access_granted"
$ws.Range("E2:E89").Value = "Opened invoice page"

$timestamps = @(
    "2025-12-25 02:10:21",
    "2025-12-25 02:10:23",
    "2025-12-25 02:10:23",
    "2025-12-25 02:10:41",
    "2025-12-25 02:10:44",
    "2025-12-25 02:10:44",
    "2025-12-25 02:10:44",
    "2025-12-25 02:10:45",
    "2025-12-25 02:10:50",
    "2025-12-25 02:10:50",
    "2025-12-25 02:10:50",
    "2025-12-25 02:10:51",
    "2025-12-25 02:10:55",
    "2025-12-25 02:10:55",
    "2025-12-25 02:11:04",
    "2025-12-25 02:11:24",
    "2025-12-25 02:11:30",
    "2025-12-25 02:11:31",
    "2025-12-25 02:11:32",
    "2025-12-25 02:11:32",
    "2025-12-25 02:11:34",
    "2025-12-25 02:11:34",
    "2025-12-25 02:11:44",
    "2025-12-25 02:11:52",
    "2025-12-25 02:11:53",
    "2025-12-25 02:11:53",
    "2025-12-25 02:12:02",
    "2025-12-25 02:12:08",
    "2025-12-25 02:12:08",
    "2025-12-25 02:12:16",
    "2025-12-25 02:12:23",
    "2025-12-25 02:12:23",
    "2025-12-25 02:12:27",
    "2025-12-25 02:13:00",
    "2025-12-25 02:13:00",
    "2025-12-25 02:13:17",
    "2025-12-25 02:13:42",
    "2025-12-25 02:13:43",
    "2025-12-25 02:13:56",
    "2025-12-25 02:14:03",
    "2025-12-25 02:14:03",
    "2025-12-25 02:14:12",
    "2025-12-25 02:14:13",
    "2025-12-25 02:14:14",
    "2025-12-25 02:14:15",
    "2025-12-25 02:14:21",
    "2025-12-25 02:14:21",
    "2025-12-25 02:14:27",
    "2025-12-25 02:14:38",
    "2025-12-25 02:14:39",
    "2025-12-25 02:14:39",
    "2025-12-25 02:14:48",
    "2025-12-25 02:14:54",
    "2025-12-25 02:14:56",
    "2025-12-25 02:14:57",
    "2025-12-25 02:15:07",
    "2025-12-25 02:15:13",
    "2025-12-25 02:15:15",
    "2025-12-25 02:15:15",
    "2025-12-25 02:15:31",
    "2025-12-25 02:15:45",
    "2025-12-25 02:15:45",
    "2025-12-25 02:15:46",
    "2025-12-25 02:15:53",
    "2025-12-25 02:16:00",
    "2025-12-25 02:16:01",
    "2025-12-25 02:16:01",
    "2025-12-25 02:16:09",
    "2025-12-25 02:16:15",
    "2025-12-25 02:16:25",
    "2025-12-25 02:16:29",
    "2025-12-25 02:16:29",
    "2025-12-25 02:16:34",
    "2025-12-25 02:16:37",
    "2025-12-25 02:16:37",
    "2025-12-25 02:16:53",
    "2025-12-25 02:16:58",
    "2025-12-25 02:16:59",
    "2025-12-25 02:17:13",
    "2025-12-25 02:17:37",
    "2025-12-25 02:17:50",
    "2025-12-25 02:18:10",
    "2025-12-25 02:18:46",
    "2025-12-25 02:21:45",
    "2025-12-25 02:21:49",
    "2025-12-25 02:21:54",
    "2025-12-25 02:23:07",
    "2025-12-25 02:23:41"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $timestamps[$i]
}
